# Logartimización de todas las variables y eliminación de categoricas
#
# This script:
#  1) Renames the predictor labels in column B (rows 3-18) removing the
#     backticks and replacing internal spaces/colons with dots.
#  2) Updates the coefficient values in column C (rows 2,3,5,6,7,9-14,16,17,18)
#     to the new, re-fit lasso coefficients.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename predictor labels (column B) ---------------------------------
$labelMap = @{
    3  = "pH.agua.suelo"
    4  = "Fósforo.Bray.II"
    5  = "Azufre.Fosfato.monocalcico"
    6  = "Acidez.Intercambiable"
    7  = "Aluminio.intercambiable"
    8  = "Calcio.intercambiable"
    9  = "Magnesio.intercambiable"
    10 = "Potasio.intercambiable"
    11 = "Sodio.intercambiable"
    12 = "capacidad.de.intercambio.cationico"
    13 = "Conductividad.electrica"
    14 = "Hierro.disponible.olsen"
    15 = "Cobre.disponible"
    16 = "Manganeso.disponible.Olsen"
    17 = "Zinc.disponible.Olsen"
    18 = "Boro.disponible"
}

foreach ($row in $labelMap.Keys) {
    $ws.Cells.Item($row, 2).Value = $labelMap[$row]
}

# --- 2. Update coefficient values (column C) -------------------------------
$valueMap = @{
    2  = 2.1454175777854183
    3  = 0.44762115175654044
    5  = 0.012413611575142378
    6  = 0.07902442609122665
    7  = 0.06512373023491637
    9  = -0.1650672261683916
    10 = 0.12292377958328944
    11 = -0.11406705779358295
    12 = -0.15889824917456358
    13 = 0.18270469928597377
    14 = 0.07593913294443458
    16 = -0.07845058366124621
    17 = 0.10460414134404898
    18 = -0.002902555094178036
}

foreach ($row in $valueMap.Keys) {
    $ws.Cells.Item($row, 3).Value = $valueMap[$row]
}
